$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.571.39"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.441.82"
$ws.Range("E3").Value = "  -2.54%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.09"
$ws.Range("E5").Value = "  -1.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.25"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.609"
$ws.Range("E7").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.440.53"
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -2.18%  "
$ws.Range("E11").Value = "  -2.75%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -3.58%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.040.41"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.94"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.572.80"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("E17").Value = "  -4.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.445.95"
$ws.Range("E18").Value = "  -2.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.95"
$ws.Range("E20").Value = "  -6.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.57"
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.82"
$ws.Range("E22").Value = "  -3.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +1.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("E24").Value = "  -0.14%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.531"
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.29"
$ws.Range("E26").Value = "  -3.18%  "
$ws.Range("E27").Value = "  -5.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  -5.37%  "
$ws.Range("E29").Value = "  -2.57%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.03"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("E33").Value = "  -7.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.20"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.20"
$ws.Range("E36").Value = "  -4.27%  "
$ws.Range("E37").Value = "  -8.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.96"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.882"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -6.80%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  -2.92%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.59"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.63"
$ws.Range("E43").Value = "  -8.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "25.72"
$ws.Range("E44").Value = "  -5.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0711"
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.92"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.696.63"
$ws.Range("E47").Value = "  -6.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.13"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0295"
$ws.Range("E49").Value = "  -4.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "325.04"
$ws.Range("E50").Value = "  -8.25%  "
$ws.Range("E51").Value = "  -5.75%  "
